$wb = $excel.ActiveWorkbook

# --- Sheet 1 ("begroting"): update weekly hours for week 17 (B40) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B40").Value = 33

# --- Sheet 2 ("Sheet1"): add spherical-coordinate dx/dtmax samples in columns C and D ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("C6").Value = 3
$ws2.Range("D6").Value = 1
$ws2.Range("C7").Formula = "=C6/SQRT(9.81)"
$ws2.Range("D7").Formula = "=D6/SQRT(9.81)"

# Move the active selection on sheet 2 to L16, matching the saved view state
$ws2.Range("L16").Select()

$wb.Save()
